# Update countries & provincias Spain
#
# 1. Refresh the "Datos actualizados" timestamp (row 1) from 08:21 to 09:38.
# 2. Refresh COVID stats for several countries (Rusia, Armenia, Moldavia,
#    Hungria, Lituania, Estonia, Letonia).
# 3. Israel/Ucrania swap places in the (descending, by "Casos totales")
#    ranking: Ucrania's updated total (68794) now exceeds Israel's
#    (68556, unchanged), so Ucrania takes row 36 and Israel moves to row 37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp -------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Julio de 2020 a las 09:38"

# --- Rusia (row 7) -----------------------------------------------------
$ws.Range("B7").Value = 834499
$ws.Range("C7").Value = 5509
$ws.Range("D7").Value = 629655
$ws.Range("E7").Value = 191042
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 129
$ws.Range("H7").Value = 13802

# --- Ucrania / Israel swap (rows 36 & 37) ------------------------------
# Row 36 becomes Ucrania with fresh numbers.
$ws.Range("A36").Value = "Ucrania"
$ws.Range("B36").Value = 68794
$ws.Range("C36").Value = 1197
$ws.Range("D36").Value = 38154
$ws.Range("E36").Value = 28967
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 23
$ws.Range("H36").Value = 1673

# Row 37 becomes Israel, carrying the numbers Israel had in the old row 36.
$ws.Range("A37").Value = "Israel"
$ws.Range("B37").Value = 68556
$ws.Range("C37").Value = 257
$ws.Range("D37").Value = 35513
$ws.Range("E37").Value = 32552
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 491

# --- Armenia (row 53) ---------------------------------------------------
$ws.Range("B53").Value = 38196
$ws.Range("C53").Value = 259
$ws.Range("D53").Value = 28366
$ws.Range("E53").Value = 9102
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 5
$ws.Range("H53").Value = 728

# --- Moldavia (row 63) ---------------------------------------------------
$ws.Range("D63").Value = 17040
$ws.Range("E63").Value = 6145
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 3
$ws.Range("H63").Value = 762

# --- Hungria (row 103) ---------------------------------------------------
$ws.Range("B103").Value = 4484
$ws.Range("C103").Value = 19
$ws.Range("D103").Value = 3346
$ws.Range("E103").Value = 542

# --- Lituania (row 126) ---------------------------------------------------
$ws.Range("B126").Value = 2062
$ws.Range("C126").Value = 19
$ws.Range("E126").Value = 339

# --- Estonia (row 127) ---------------------------------------------------
$ws.Range("B127").Value = 2051
$ws.Range("C127").Value = 9
$ws.Range("E127").Value = 56

# --- Letonia (row 140) ---------------------------------------------------
$ws.Range("B140").Value = 1228
$ws.Range("C140").Value = 4
$ws.Range("E140").Value = 145
